# "progress of the day" edit for the dates.xlsx workbook
#
# The underlying data had a duplicate "Marketing action 3" row appended at
# the end of each Course Year group (rows 5, 9 and 13 in the original
# 1-based row numbering). Cleaning the dataset means dropping those three
# duplicate rows, which shifts everything below them up.
#
# After the cleanup, an AutoFilter is turned on over the now-smaller table
# (A1:F10), which is also what causes Excel to persist the hidden
# "_xlnm._FilterDatabase" workbook-scoped (sheet-local) defined name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the duplicate rows bottom-up so the remaining row numbers don't
# shift out from under us while we work.
$ws.Rows("13").Delete()
$ws.Rows("9").Delete()
$ws.Rows("5").Delete()

# Turn on AutoFilter for the cleaned-up table.
$ws.Range("A1:F10").AutoFilter()

# Excel records the AutoFilter range as a hidden, sheet-scoped defined name.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=dates!`$A`$1:`$F`$10")
$filterName.Visible = $false

# Leave the cursor where the author left it when they saved.
$ws.Range("D13").Select()
